$wb = $excel.ActiveWorkbook

# "Generate Report for Archive": the localization status moves from
# "Ready for handoff" to "In Translation" for the e2e markdown file.
# This value appears in the Status column on all three sheets:
#   - Overview: columns "zh-cn" (E2) and "de-de" (F2)
#   - zh-cn:    column "Status" (C2)
#   - de-de:    column "Status" (C2)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsZhCn.Range("C2").Value = "In Translation"
$wsDeDe.Range("C2").Value = "In Translation"

# The Status columns were re-sized (narrower, since "In Translation" is
# shorter than "Ready for handoff"). Apply the equivalent column widths.
$wsOverview.Range("E1").ColumnWidth = 12.5
$wsOverview.Range("F1").ColumnWidth = 12.5
$wsZhCn.Range("C1").ColumnWidth = 12.5
$wsDeDe.Range("C1").ColumnWidth = 12.5
